# Insert a new data row at row 26 (pushing existing rows 26-57 down to 27-58)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 44589
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 100112022
$ws.Cells.Item(26, 7).Value = "Arveja Verde"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 120
$ws.Cells.Item(26, 11).Value = 23000
$ws.Cells.Item(26, 12).Value = 24000
$ws.Cells.Item(26, 13).Value = 23500
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(26, 16).Value = 940
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
